$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2928.1667
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 3313.8
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 9941.400000000001
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -10277.4
# Row 69
$ws.Range("H69").Value = 176666.67
$ws.Range("I69").Value = 10000
$ws.Range("K69").Value = 30000
$ws.Range("M69").Value = -29126
# Row 72
$ws.Range("H72").Value = 176666.67
$ws.Range("I72").Value = 10000
$ws.Range("K72").Value = 90000
$ws.Range("M72").Value = -85632
# Row 80
$ws.Range("H80").Value = 3783.3333
$ws.Range("I80").Value = 1200
$ws.Range("J80").Value = 6366.6665
$ws.Range("K80").Value = 3600
$ws.Range("L80").Value = 19099.9995
$ws.Range("M80").Value = -2602
$ws.Range("N80").Value = -21095.9995
# Row 83
$ws.Range("H83").Value = 3783.3333
$ws.Range("I83").Value = 1200
$ws.Range("J83").Value = 6366.6665
$ws.Range("K83").Value = 10800
$ws.Range("L83").Value = 57299.9985
$ws.Range("M83").Value = -5808
$ws.Range("N83").Value = -67283.9985
# Row 87
$ws.Range("H87").Value = 99353
$ws.Range("J87").Value = 99353
$ws.Range("L87").Value = 99353
$ws.Range("N87").Value = -101849
# Row 90
$ws.Range("H90").Value = 99353
$ws.Range("J90").Value = 99353
$ws.Range("L90").Value = 298059
$ws.Range("N90").Value = -310539
# Row 106
$ws.Range("H106").Value = 3999
$ws.Range("I106").Value = 3999
$ws.Range("K106").Value = 3999
$ws.Range("M106").Value = -3368
# Row 116
$ws.Range("H116").Value = 140714.28
$ws.Range("I116").Value = 140714.28
$ws.Range("K116").Value = 140714.28
$ws.Range("M116").Value = -137272.28
# Row 141
$ws.Range("H141").Value = 7998.6
$ws.Range("I141").Value = 7998.6
$ws.Range("K141").Value = 23995.8
$ws.Range("M141").Value = -18815.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3911.205
$ws.Range("I32").Value = 3816.7896
$ws.Range("K32").Value = 3816.7896
$ws.Range("M32").Value = -3529.7896
# Row 97
$ws.Range("H97").Value = 1917.8
$ws.Range("I97").Value = 1045
$ws.Range("K97").Value = 1045
$ws.Range("M97").Value = -549
# Row 132
$ws.Range("H132").Value = 2853.3125
$ws.Range("I132").Value = 2138.0833
$ws.Range("K132").Value = 6414.249899999999
$ws.Range("M132").Value = -3884.249899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 3421.8
$ws.Range("I94").Value = 3536.3333
$ws.Range("J94").Value = 3250
$ws.Range("K94").Value = 3536.3333
$ws.Range("L94").Value = 3250
$ws.Range("M94").Value = -3085.3333
$ws.Range("N94").Value = -4152
# Row 99
$ws.Range("H99").Value = 649.5
$ws.Range("I99").Value = 649.5
$ws.Range("K99").Value = 649.5
$ws.Range("M99").Value = 848.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 5399.75
$ws.Range("I16").Value = 3299
$ws.Range("J16").Value = 6100
$ws.Range("K16").Value = 3299
$ws.Range("L16").Value = 6100
$ws.Range("M16").Value = -3012
$ws.Range("N16").Value = -6674
# Row 31
$ws.Range("H31").Value = 1102.1333
$ws.Range("I31").Value = 857.4545000000001
$ws.Range("J31").Value = 1775
$ws.Range("K31").Value = 857.4545000000001
$ws.Range("L31").Value = 1775
$ws.Range("M31").Value = -562.4545000000001
$ws.Range("N31").Value = -2365
# Row 34
$ws.Range("H34").Value = 1102.1333
$ws.Range("I34").Value = 857.4545000000001
$ws.Range("J34").Value = 1775
$ws.Range("K34").Value = 857.4545000000001
$ws.Range("L34").Value = 1775
$ws.Range("M34").Value = -655.4545000000001
$ws.Range("N34").Value = -2179
# Row 95
$ws.Range("H95").Value = 36424.5
$ws.Range("J95").Value = 36424.5
$ws.Range("L95").Value = 36424.5
$ws.Range("N95").Value = -41916.5
# Row 113
$ws.Range("H113").Value = 5399.75
$ws.Range("I113").Value = 3299
$ws.Range("J113").Value = 6100
$ws.Range("K113").Value = 3299
$ws.Range("L113").Value = 6100
$ws.Range("M113").Value = -1129
$ws.Range("N113").Value = -10440
# Row 132
$ws.Range("H132").Value = 3445.6667
$ws.Range("I132").Value = 2352.1667
$ws.Range("K132").Value = 7056.500100000001
$ws.Range("M132").Value = -4526.500100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 126.833336
$ws.Range("I2").Value = 97.333336
$ws.Range("K2").Value = 584.000016
$ws.Range("M2").Value = -471.000016
# Row 26
$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 300
$ws.Range("K26").Value = 900
$ws.Range("M26").Value = -612
# Row 37
$ws.Range("H37").Value = 75000
$ws.Range("J37").Value = 75000
$ws.Range("L37").Value = 225000
$ws.Range("N37").Value = -225224
# Row 50
$ws.Range("H50").Value = 500
$ws.Range("I50").Value = 500
$ws.Range("K50").Value = 1500
$ws.Range("M50").Value = -1019
# Row 53
$ws.Range("H53").Value = 500
$ws.Range("I53").Value = 500
$ws.Range("K53").Value = 1500
$ws.Range("M53").Value = -1019

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 114.72727
$ws.Range("I2").Value = 132.75
$ws.Range("K2").Value = 132.75
$ws.Range("M2").Value = -19.75
# Row 80
$ws.Range("H80").Value = 3099
$ws.Range("I80").Value = 3099
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3099
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2101
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 3099
$ws.Range("I83").Value = 3099
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15495
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10503
$ws.Range("N83").ClearContents()
# Row 101
$ws.Range("H101").Value = 89999
$ws.Range("J101").Value = 89999
$ws.Range("L101").Value = 89999
$ws.Range("N101").Value = -96489

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 42
$ws.Range("H42").Value = 65000
$ws.Range("J42").Value = 65000
$ws.Range("L42").Value = 65000
$ws.Range("N42").Value = -66126
# Row 49
$ws.Range("H49").Value = 65000
$ws.Range("J49").Value = 65000
$ws.Range("L49").Value = 65000
$ws.Range("N49").Value = -65294
# Row 55
$ws.Range("H55").Value = 1000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3570.5715
$ws.Range("I132").Value = 1932.6666
$ws.Range("K132").Value = 5797.9998
$ws.Range("M132").Value = -3267.9998
